$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "[[0.44889334]`n [0.441223  ]]"
$ws.Range("C2").Value = "[[0.20702952]`n [0.17418403]]"
$ws.Range("D2").Value = 6000
